# ------------------------------------------------------------------
# Applies the "correccion beneficios al negocio" edit:
#   1) Merge two adjacent runs ("... para permitir que puedan," and
#      " ante distintos inconvenientes, ... oportuno.") into a single
#      run with the combined text (same run formatting on both sides).
#   2) Replace the "Puede ser usado de plataforma ..." sentence with
#      the new "Tiene el potencial de generar nuevos negocios ..."
#      text, split across six runs (all sharing the original run's
#      formatting), matching the target XML structure.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Part 1 -----------------------------------------------------------
# The two original runs already carry identical run formatting, so a
# single in-place replacement (wdReplaceOne) that reproduces the exact
# concatenated text causes the engine to fold the second run into the
# first, leaving one run behind - exactly what the diff shows.
$part1Old = " ante distintos inconvenientes, tomar las acciones preventivas o correctivas necesarias en el momento oportuno."
$rng1 = $d.Content
$rng1.Find.Execute($part1Old, $true, $false, $false, $false, $false, $true, 1, $false, $part1Old, 1) | Out-Null

# --- Part 2 -----------------------------------------------------------
$part2Old = "Puede ser usado de plataforma para múltiples servicios adicionales relacionados, como recomendaciones personalizadas sobre el uso de agua, o integraciones con sistemas de control hogareños para el encendido y apagado de la bomba de suministro."
$part2New = "Tiene el potencial de generar nuevos negocios basados en el consumo de métricas, como la generación de recomendaciones personalizadas sobre el uso de agua, o integraciones con sistemas de control hogareños para la gestión automatizada del tanque."

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($part2Old, $true, $false, $false, $false, $false, $true, 1, $false, $part2New, 1)

if ($found2) {
    $base = $rng2.Start

    # Character offsets (relative to $base) delimiting the six runs that
    # the new sentence must be split into, all sharing the same rPr as
    # the original single run.
    $segments = @(
        "Tiene el potencial de generar nuevos negocios basados en el consumo de métricas, ",
        "como ",
        "la generación de ",
        "recomendaciones personalizadas sobre el uso de agua, o integraciones con sistemas de control hogareños para",
        " la gestión automatizada del tanque",
        "."
    )

    $pos = 0
    foreach ($seg in $segments) {
        $segStart = $base + $pos
        $segEnd = $segStart + $seg.Length
        $segRange = $d.Range($segStart, $segEnd)
        # Toggling a character property on/off forces the run to be
        # split at these boundaries without altering the final
        # formatting (it ends up identical to the surrounding runs).
        $segRange.Font.Bold = 1
        $segRange.Font.Bold = 0
        $pos = $pos + $seg.Length
    }
}
